# Updates leve-profit market data cells (columns H-N) across the
# per-job worksheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR), matching the
# latest scheduled-runner market pull. CUL has no changes this run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 5097.875
$ws.Range("I4").Value = 5097.875
$ws.Range("K4").Value = 5097.875
$ws.Range("M4").Value = -4983.875
$ws.Range("H17").Value = 7695146
$ws.Range("J17").Value = 7695146
$ws.Range("L17").Value = 23085438
$ws.Range("N17").Value = -23085774
$ws.Range("H33").Value = 10673.226
$ws.Range("I33").Value = 11767.107
$ws.Range("J33").Value = 463.66666
$ws.Range("K33").Value = 11767.107
$ws.Range("L33").Value = 463.66666
$ws.Range("M33").Value = -11538.107
$ws.Range("N33").Value = -921.66666
$ws.Range("H48").Value = 2636.2727
$ws.Range("J48").Value = 2636.2727
$ws.Range("L48").Value = 7908.8181
$ws.Range("N48").Value = -8492.8181
$ws.Range("H56").Value = 2636.2727
$ws.Range("J56").Value = 2636.2727
$ws.Range("L56").Value = 7908.8181
$ws.Range("N56").Value = -8976.8181
$ws.Range("H81").Value = 74979.664
$ws.Range("J81").Value = 74979.664
$ws.Range("L81").Value = 74979.664
$ws.Range("N81").Value = -76975.664
$ws.Range("H84").Value = 74979.664
$ws.Range("J84").Value = 74979.664
$ws.Range("L84").Value = 224938.992
$ws.Range("N84").Value = -234922.992
$ws.Range("H86").Value = 3278.4443
$ws.Range("I86").Value = 3259.3845
$ws.Range("K86").Value = 3259.3845
$ws.Range("M86").Value = -2136.3845
$ws.Range("H89").Value = 3278.4443
$ws.Range("I89").Value = 3259.3845
$ws.Range("K89").Value = 16296.9225
$ws.Range("M89").Value = -10680.9225
$ws.Range("H132").Value = 2123.0344
$ws.Range("I132").Value = 1663.16
$ws.Range("J132").Value = 4997.25
$ws.Range("K132").Value = 4989.48
$ws.Range("L132").Value = 14991.75
$ws.Range("M132").Value = -2459.48
$ws.Range("N132").Value = -20051.75
$ws.Range("H133").Value = 99374.5
$ws.Range("J133").Value = 99374.5
$ws.Range("L133").Value = 99374.5
$ws.Range("N133").Value = -109494.5
$ws.Range("H138").Value = 3803.95
$ws.Range("I138").Value = 1719
$ws.Range("K138").Value = 5157
$ws.Range("M138").Value = -17

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4478
$ws.Range("I2").Value = 4013.3333
$ws.Range("K2").Value = 4013.3333
$ws.Range("M2").Value = -3900.3333
$ws.Range("H32").Value = 4602.772
$ws.Range("I32").Value = 4042.8909
$ws.Range("J32").Value = 19999.5
$ws.Range("K32").Value = 4042.8909
$ws.Range("L32").Value = 19999.5
$ws.Range("M32").Value = -3755.8909
$ws.Range("N32").Value = -20573.5
$ws.Range("H44").Value = 31000
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 31000
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 31000
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -31976
$ws.Range("H74").Value = 4597.3125
$ws.Range("I74").Value = 3111.2144
$ws.Range("K74").Value = 3111.2144
$ws.Range("M74").Value = -2237.2144
$ws.Range("H77").Value = 4597.3125
$ws.Range("I77").Value = 3111.2144
$ws.Range("K77").Value = 15556.072
$ws.Range("M77").Value = -11188.072
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H97").Value = 1321.425
$ws.Range("I97").Value = 1336.5161
$ws.Range("J97").Value = 1269.4445
$ws.Range("K97").Value = 1336.5161
$ws.Range("L97").Value = 1269.4445
$ws.Range("M97").Value = -840.5161000000001
$ws.Range("N97").Value = -2261.4445
$ws.Range("H102").Value = 6332.25
$ws.Range("I102").Value = 3499.75
$ws.Range("K102").Value = 3499.75
$ws.Range("M102").Value = -1877.75
$ws.Range("H116").Value = 4478
$ws.Range("I116").Value = 4013.3333
$ws.Range("K116").Value = 4013.3333
$ws.Range("M116").Value = -1719.3333
$ws.Range("H122").Value = 1468.92
$ws.Range("J122").Value = 1599.5
$ws.Range("L122").Value = 4798.5
$ws.Range("N122").Value = -9698.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4478
$ws.Range("I3").Value = 4013.3333
$ws.Range("K3").Value = 4013.3333
$ws.Range("M3").Value = -3899.3333
$ws.Range("H94").Value = 1913.3334
$ws.Range("I94").Value = 1777.5
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 1777.5
$ws.Range("L94").Value = 3000
$ws.Range("M94").Value = -1326.5
$ws.Range("N94").Value = -3902
$ws.Range("H96").Value = 20284.6
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H97").Value = 12006.5
$ws.Range("I97").Value = 7093.143
$ws.Range("K97").Value = 7093.143
$ws.Range("M97").Value = -6102.143
$ws.Range("H107").Value = 3161.4
$ws.Range("I107").Value = 2900.125
$ws.Range("K107").Value = 2900.125
$ws.Range("M107").Value = -980.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 384.22223
$ws.Range("I22").Value = 407.25
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 407.25
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = -57.25
$ws.Range("N22").Value = -900
$ws.Range("H31").Value = 66673150
$ws.Range("I31").Value = 100003480
$ws.Range("K31").Value = 100003480
$ws.Range("M31").Value = -100003185
$ws.Range("H34").Value = 66673150
$ws.Range("I34").Value = 100003480
$ws.Range("K34").Value = 100003480
$ws.Range("M34").Value = -100003278
$ws.Range("H58").Value = 9812.956
$ws.Range("J58").Value = 11997.333
$ws.Range("L58").Value = 11997.333
$ws.Range("N58").Value = -12403.333
$ws.Range("H69").Value = 26538.666
$ws.Range("J69").Value = 29871.2
$ws.Range("L69").Value = 29871.2
$ws.Range("N69").Value = -31369.2
$ws.Range("H72").Value = 26538.666
$ws.Range("J72").Value = 29871.2
$ws.Range("L72").Value = 89613.60000000001
$ws.Range("N72").Value = -97101.60000000001
$ws.Range("H93").Value = 20198.8
$ws.Range("I93").Value = 10333
$ws.Range("J93").Value = 34997.5
$ws.Range("K93").Value = 10333
$ws.Range("L93").Value = 34997.5
$ws.Range("M93").Value = -8461
$ws.Range("N93").Value = -38741.5
$ws.Range("H112").Value = 46287.855
$ws.Range("J112").Value = 46287.855
$ws.Range("L112").Value = 46287.855
$ws.Range("N112").Value = -49241.855
$ws.Range("H133").Value = 62713
$ws.Range("J133").Value = 62713
$ws.Range("L133").Value = 62713
$ws.Range("N133").Value = -67773
$ws.Range("H136").Value = 9812.956
$ws.Range("J136").Value = 11997.333
$ws.Range("L136").Value = 35991.999
$ws.Range("N136").Value = -41091.999
$ws.Range("H141").Value = 212262.1
$ws.Range("J141").Value = 212262.1
$ws.Range("L141").Value = 212262.1
$ws.Range("N141").Value = -222622.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 44945
$ws.Range("J25").Value = 44945
$ws.Range("L25").Value = 44945
$ws.Range("N25").Value = -46003
$ws.Range("H80").Value = 3155.889
$ws.Range("I80").Value = 3691.125
$ws.Range("J80").Value = 2727.7
$ws.Range("K80").Value = 3691.125
$ws.Range("L80").Value = 2727.7
$ws.Range("M80").Value = -2693.125
$ws.Range("N80").Value = -4723.7
$ws.Range("H83").Value = 3155.889
$ws.Range("I83").Value = 3691.125
$ws.Range("J83").Value = 2727.7
$ws.Range("K83").Value = 18455.625
$ws.Range("L83").Value = 13638.5
$ws.Range("M83").Value = -13463.625
$ws.Range("N83").Value = -23622.5
$ws.Range("H97").Value = 1304.6
$ws.Range("I97").Value = 1040.7778
$ws.Range("K97").Value = 1040.7778
$ws.Range("M97").Value = -544.7778000000001
$ws.Range("H100").Value = 105995
$ws.Range("J100").Value = 105995
$ws.Range("L100").Value = 105995
$ws.Range("N100").Value = -108159

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 115414.336
$ws.Range("I61").Value = 129466.25
$ws.Range("K61").Value = 129466.25
$ws.Range("M61").Value = -129264.25
$ws.Range("H68").Value = 3760.8
$ws.Range("I68").Value = 2268.3333
$ws.Range("K68").Value = 2268.3333
$ws.Range("M68").Value = -1519.3333
$ws.Range("H71").Value = 3760.8
$ws.Range("I71").Value = 2268.3333
$ws.Range("K71").Value = 11341.6665
$ws.Range("M71").Value = -7597.666499999999
$ws.Range("H113").Value = 115414.336
$ws.Range("I113").Value = 129466.25
$ws.Range("K113").Value = 129466.25
$ws.Range("M113").Value = -127296.25
$ws.Range("H136").Value = 7107.2915
$ws.Range("I136").Value = 7069.0435
$ws.Range("K136").Value = 21207.1305
$ws.Range("M136").Value = -18657.1305

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 18833.666
$ws.Range("I41").Value = 18200
$ws.Range("J41").Value = 18912.875
$ws.Range("K41").Value = 18200
$ws.Range("L41").Value = 18912.875
$ws.Range("M41").Value = -17810
$ws.Range("N41").Value = -19692.875
$ws.Range("H99").Value = 39833.25
$ws.Range("J99").Value = 39859.668
$ws.Range("L99").Value = 39859.668
$ws.Range("N99").Value = -45849.668
